# "Update countries & provincias Spain" - refresh the COVID figures in the
# "Pais" sheet and re-sort a handful of rows whose totals crossed a
# neighbouring country's total (same effect as the source re-generating the
# sheet from freshly-sorted data: the row keeps its position, but the
# country label + stats that land on it change).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner (row 1)
$ws.Range("A1").Value = "Datos actualizados a 8 de Octubre de 2020 a las 06:20"

# India (row 5) - refreshed totals, same rank
$ws.Range("B5").Value = 6835655
$ws.Range("C5").Value = 2667
$ws.Range("D5").Value = 5827704
$ws.Range("E5").Value = 902397

# Honduras overtakes Etiopia (rows 54-55 swap place + figures)
$ws.Range("A54").Value = "Honduras"
$ws.Range("B54").Value = 81016
$ws.Range("C54").Value = 354
$ws.Range("D54").Value = 30590
$ws.Range("E54").Value = 47960
$ws.Range("G54").Value = 19
$ws.Range("H54").Value = 2466
$ws.Range("A55").Value = "Etiopia"
$ws.Range("B55").Value = 80895
$ws.Range("D55").Value = 35670
$ws.Range("E55").Value = 43970
$ws.Range("H55").Value = 1255

# Tailandia (row 142) - refreshed totals, same rank
$ws.Range("B142").Value = 3622
$ws.Range("C142").Value = 7
$ws.Range("D142").Value = 3439
$ws.Range("E142").Value = 124

# Belice overtakes Sierra Leona and Letonia (rows 154-156 shift down one + figures)
$ws.Range("A154").Value = "Belice"
$ws.Range("B154").Value = 2310
$ws.Range("C154").Value = 67
$ws.Range("D154").Value = 1427
$ws.Range("E154").Value = 849
$ws.Range("H154").Value = 34
$ws.Range("A155").Value = "Sierra Leona"
$ws.Range("B155").Value = 2287
$ws.Range("D155").Value = 1716
$ws.Range("E155").Value = 499
$ws.Range("H155").Value = 72
$ws.Range("A156").Value = "Letonia"
$ws.Range("B156").Value = 2261
$ws.Range("D156").Value = 1322
$ws.Range("E156").Value = 899
$ws.Range("H156").Value = 40

# Butan (row 187) - refreshed totals, same rank
$ws.Range("B187").Value = 304
$ws.Range("C187").Value = 4
$ws.Range("D187").Value = 252
$ws.Range("E187").Value = 52

# Santa Lucia overtakes Nueva Caledonia (rows 207-208 swap place, figures tied)
$ws.Range("A207").Value = "Santa Lucia"
$ws.Range("A208").Value = "Nueva Caledonia"

# Islas Malvinas overtakes Montserrat (rows 215-216 swap place + figures)
$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0
$ws.Range("A216").Value = "Montserrat"
$ws.Range("D216").Value = 12
$ws.Range("H216").Value = 1
